$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The single "T7" category used on rows 3-6 is renamed to "t7d" now that
# more time points have been added to the experiment.
for ($r = 3; $r -le 6; $r++) {
    $ws.Range("C$r").Value = "t7d"
}

# New data rows 7-15: three more time-point groups (t24h, t0h, t6h), each
# with three replicate weight measurements, sharing the "vulgaris" morfotipo
# and a new "fecha de corte" of 2025-09-19 (serial 45919).
$newRows = @(
    @{ Row=7;  A=5;  B=51; C="t24h" },
    @{ Row=8;  A=6;  B=36; C="t24h" },
    @{ Row=9;  A=7;  B=13; C="t24h" },
    @{ Row=10; A=8;  B=68; C="t0h"  },
    @{ Row=11; A=9;  B=46; C="t0h"  },
    @{ Row=12; A=10; B=41; C="t0h"  },
    @{ Row=13; A=11; B=15; C="t6h"  },
    @{ Row=14; A=12; B=19; C="t6h"  },
    @{ Row=15; A=13; B=43; C="t6h"  }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "vulgaris"
    $ws.Cells.Item($r, 5).Value = 45919
}

# Copy the existing date formatting from E3:E6 onto the new date cells so
# they reuse the same style (numFmtId 14) instead of creating a new one.
$ws.Range("E3:E6").Copy()
$ws.Range("E7:E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Summary rows: average and sample standard deviation of the "peso (g)"
# column across all replicates.
$ws.Range("B17").Formula = "=STDEV.S(B3:B15)"
$ws.Range("B17").NumberFormat = "0.000"

$ws.Range("B16").Formula = "=AVERAGE(B3:B15)"
$ws.Range("B16").NumberFormat = "0.00"

$ws.Range("B17").Select()
